# Daily attendance processing - 2025-09-29 12:18:55
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PercentText {
    # Writing a "NN.N%"-looking string via .Value auto-converts to a numeric
    # percentage (real Excel behaviour). Force Text format first so it lands
    # as literal text, then paste the *formatting* back from a same-style
    # neighbour cell so the original cell style index is preserved.
    param($addr, $value, $fmtSourceAddr)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $ws.Range($fmtSourceAddr).Copy()
    $r.PasteSpecial(-4122)
}

# --- Row 5 (Year2/A1/PHYSIOLOGY) ---
$ws.Range("G5").Value = "youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("H5").Value = "0/217"

# --- Row 6 (Year2/A1/POS): Recorded -> Pending, restyle from a Pending row ---
$ws.Range("A2").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)
$ws.Range("H6").Value = "0/217"
$ws.Range("I6").Value = "Pending"

# --- Class Statistics block (K/L column) ---
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 30
Set-PercentText "L9" "27.3%" "K9"

# --- Row 10 (Year2/A2/PHYSIOLOGY) ---
$ws.Range("G10").Value = "youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
Set-PercentText "L10" "0.0%" "K9"

# --- Row 12 (Year2/A2/POS): Recorded -> Pending, restyle from a Pending row ---
$ws.Range("A2").Copy()
$ws.Range("A12:I12").PasteSpecial(-4122)
$ws.Range("I12").Value = "Pending"

# --- Group Statistics block (rows 15,16,20,21,22) ---
$ws.Range("O15").Value = 3
$ws.Range("Q15").Value = 2
Set-PercentText "R15" "60.0%" "K9"
Set-PercentText "S15" "0.0%" "K9"

$ws.Range("O16").Value = 3
$ws.Range("Q16").Value = 3
Set-PercentText "R16" "50.0%" "K9"

$ws.Range("O20").Value = 1
$ws.Range("Q20").Value = 3
Set-PercentText "R20" "25.0%" "K9"

$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 3

$ws.Range("O22").Value = 2
$ws.Range("P22").Value = 0
Set-PercentText "R22" "40.0%" "K9"

# --- Row 27 (Year2/B1/HISTOLOGY): add recorder ---
$ws.Range("G27").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# --- Row 32 (Year2/B2/HISTOLOGY): add recorder ---
$ws.Range("G32").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

# --- Row 35 (Year2/B2/POS): Recorded -> Pending, restyle from a Pending row ---
$ws.Range("A2").Copy()
$ws.Range("A35:I35").PasteSpecial(-4122)
$ws.Range("I35").Value = "Pending"

# --- Row 39 (Year2/B3/PHYSIOLOGY): Not Recorded -> Recorded, restyle from a Recorded row ---
$ws.Range("A3").Copy()
$ws.Range("A39:I39").PasteSpecial(-4122)
$ws.Range("G39").Value = "aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("I39").Value = "Recorded"

# --- Row 40 (Year2/B3/POS): Recorded -> Pending, restyle from a Pending row ---
$ws.Range("A2").Copy()
$ws.Range("A40:I40").PasteSpecial(-4122)
$ws.Range("I40").Value = "Pending"

# --- Row 44 (Year2/B4/PHYSIOLOGY): Not Recorded -> Recorded, restyle from a Recorded row ---
$ws.Range("A3").Copy()
$ws.Range("A44:I44").PasteSpecial(-4122)
$ws.Range("G44").Value = "aya.hanafy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("I44").Value = "Recorded"
